$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: relabel the Cone section's first results row with the new cone description
$ws.Range("G11").Value = "Polystyrene cone in water (Radius 1, Height 2 micro m)"

# New row 15: final results row for the cone section
$ws.Range("A15").Value = 15
$ws.Range("B15").Value = 1.0640000000000001
$ws.Range("C15").Value = 1.1859519223999999
$ws.Range("D15").Value = 2
$ws.Range("G15").Value = "Final Results for Polystyrene cone in water (Radius 1, Height 2 micro m)"

# Row 2: add "Just considering F_z" note next to the Sphere header
$ws.Range("B2").Value = "Just considering F_z"

# Row 10: add "Considering F_x,F_y,F_z" note next to the Cone header
$ws.Range("B10").Value = "Considering F_x,F_y,F_z"

# Row 12/13: add labels for the default dpl test rows under Cone
$ws.Range("G12").Value = "Default 15 dpl tests"
$ws.Range("G13").Value = "Default 30 dpl tests"

# Update selection to match the new active cell in the diff
$ws.Range("E15").Select()
